$wb = $excel.ActiveWorkbook

# --- Rename existing "Acc_Upfront" sheet to "Acc_Upfront2" ---
$ws2 = $wb.Worksheets.Item("Acc_Upfront")
$ws2.Name = "Acc_Upfront2"

# --- Add a new sheet "Acc_Upfront1" right after it ---
$ws1 = $wb.Worksheets.Add($null, $ws2)
$ws1.Name = "Acc_Upfront1"

# --- Populate the new "Acc_Upfront1" sheet with the header row ---
$ws1.Cells.Item(1,1).Value = "Entry ID"
$ws1.Cells.Item(1,2).Value = "Office"
$ws1.Cells.Item(1,3).Value = "Transaction Date"
$ws1.Cells.Item(1,4).Value = "Transaction ID"
$ws1.Cells.Item(1,5).Value = "Type"
$ws1.Cells.Item(1,6).Value = "Created By"
$ws1.Cells.Item(1,7).Value = "Account"
$ws1.Cells.Item(1,8).Value = "Debit"
$ws1.Cells.Item(1,9).Value = "Credit"

# Row 2 (previously row 5 on Acc_Upfront2)
$ws1.Cells.Item(2,1).Value = 76
$ws1.Cells.Item(2,2).Value = "Head Office"
$ws1.Cells.Item(2,3).Value = 42005
$ws1.Cells.Item(2,3).NumberFormat = "d-mmm-yy"
$ws1.Cells.Item(2,4).Value = 3
$ws1.Cells.Item(2,5).Value = "ASSET"
$ws1.Cells.Item(2,6).Value = "mifos"
$ws1.Cells.Item(2,7).Value = "Fees Receivable(4)"
$ws1.Cells.Item(2,8).Value = "'$ 100"
$ws1.Cells.Item(2,8).Style = "Normal"

# Row 3 (previously row 6 on Acc_Upfront2)
$ws1.Cells.Item(3,1).Value = 77
$ws1.Cells.Item(3,2).Value = "Head Office"
$ws1.Cells.Item(3,3).Value = 42005
$ws1.Cells.Item(3,3).NumberFormat = "d-mmm-yy"
$ws1.Cells.Item(3,4).Value = 3
$ws1.Cells.Item(3,5).Value = "INCOME"
$ws1.Cells.Item(3,6).Value = "mifos"
$ws1.Cells.Item(3,7).Value = "Income from fees(8)"
$ws1.Cells.Item(3,9).Value = "'$ 100"
$ws1.Cells.Item(3,9).Style = "Normal"

# --- Remove the now-duplicated rows 5:6 from "Acc_Upfront2" ---
$ws2.Range("A5:I6").EntireRow.Delete()

# --- Fix up the sheet that used to be tabSelected (Acc_Repayment) ---
$wsRepay = $wb.Worksheets.Item("Acc_Repayment")

# --- View/selection bookkeeping to match the final, saved UI state ---
[void]$ws1.Activate()
[void]$ws1.Range("I23").Select()

[void]$ws2.Activate()
[void]$ws2.Range("G15").Select()
